# Console log statement: add new player "Farah Barakat" to the All Players
# roster and refresh all the sheets that are derived from it (Selected
# Team, Randomly Selected Players, Count Players by Position, Players
# Sorted by APT, Player with Lowest AVG).

function Test-NumericLike($val) {
    # Values like "23.0" / "1.0" need to stay TEXT (they are AVG strings),
    # not be auto-converted to numbers by Excel when assigned.
    if ($val -is [string]) {
        return $val -match '^-?\d+(\.\d+)?$'
    }
    return $false
}

function Set-RowValues($ws, $r, $vals) {
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $val = $vals[$i]
        $col = $i + 1
        if (Test-NumericLike $val) {
            # Leading apostrophe forces Excel to keep it as text (quote-prefix),
            # matching how the existing AVG column values are stored as text.
            $ws.Cells.Item($r, $col).Value = "'" + $val
        } else {
            $ws.Cells.Item($r, $col).Value = $val
        }
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "All Players" — append the newly registered player.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Players")
Set-RowValues $wsAll 25 @(1722160739577, "Farah", "Barakat", 12, 34, "Northern Ireland", "23.0", "Defender")

# ---------------------------------------------------------------------
# 2. "Selected Team" — append the newly (re)computed team roster rows.
# ---------------------------------------------------------------------
$wsTeam = $wb.Worksheets.Item("Selected Team")
Set-RowValues $wsTeam 5  @(4, "Jordan", "Robinson", 45, 89, "Wales", "67.0", "Attacker")
Set-RowValues $wsTeam 6  @(5, "Steven", "Walker", 88, 87, "Wales", "87.5", "Midfielder")
Set-RowValues $wsTeam 7  @(7, "Rashid", "Bhatti", 90, 86, "England", "88.0", "Midfielder")
Set-RowValues $wsTeam 8  @(8, "Thomas", "Taylor", 97, 85, "England", "91.0", "Defender")
Set-RowValues $wsTeam 9  @(9, "Theo", "Dolan", 87, 82, "Scotland", "84.5", "Attacker")
Set-RowValues $wsTeam 10 @(6, "Alfie", "Loy", 85, 79, "Wales", "82.0", "Attacker")

# ---------------------------------------------------------------------
# 3. "Randomly Selected Players" — new random draw.
# ---------------------------------------------------------------------
$wsRand = $wb.Worksheets.Item("Randomly Selected Players")
Set-RowValues $wsRand 2 @(14, "William", "Adams", 78, 78, "England", "78.0", "Midfielder")
Set-RowValues $wsRand 3 @(13, "Isaac", "Johnson", 76, 77, "England", "76.5", "Defender")

# ---------------------------------------------------------------------
# 4. "Count Players by Position" — Defender count goes up by 1.
# ---------------------------------------------------------------------
$wsCount = $wb.Worksheets.Item("Count Players by Position")
$wsCount.Cells.Item(2, 2).Value = 8

# ---------------------------------------------------------------------
# 5. "Players Sorted by APT" — resorted, with the new player and a
#    reshuffled tail of the tied (APT=1) entries.
# ---------------------------------------------------------------------
$wsSorted = $wb.Worksheets.Item("Players Sorted by APT")
Set-RowValues $wsSorted 21 @(1722160739577, "Farah", "Barakat", 12, 34, "Northern Ireland", "23.0", "Defender")
Set-RowValues $wsSorted 22 @(1722093683411, "nana", "b", 1, 1, "England", "1.0", "Defender")
Set-RowValues $wsSorted 23 @(1722095244591, "rana", "barakat", 1, 1, "England", "1.0", "Defender")
Set-RowValues $wsSorted 24 @(1722093666927, "r", "m", 1, 1, "England", "1.0", "Defender")
Set-RowValues $wsSorted 25 @(1722094926709, "rasha", "barakattt", 1, 1, "England", "1.0", "Defender")

# ---------------------------------------------------------------------
# 6. "Player with Lowest AVG" — the tie-break among AVG=1.0 picks a new
#    player now that the roster changed.
# ---------------------------------------------------------------------
$wsLowAvg = $wb.Worksheets.Item("Player with Lowest AVG")
Set-RowValues $wsLowAvg 2 @(1722093683411, "nana", "b", 1, 1, "England", "1.0", "Defender")

# "Player with Highest APT" is unchanged by this edit.
